# Grading workbook update:
#  - Row 20 (question 12, "For successfully...findAllBrands()" row) grading
#    comment changes from "(-5)For incorrect logic" to a new, more specific
#    comment, and the awarded points for that row drop from 5 to 1.
#  - Selection/active cell moves to F20 to reflect where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the points awarded for row 20 (column E) from 5 to 1.
$ws.Range("E20").Value = 1

# Update the grading comment for row 20 (column F) with the new remark.
$ws.Range("F20").Value = "(-9) For the logic inside the method is totally incorrect."

# Reflect the new active cell/selection used while making this edit.
$ws.Range("F20").Select()
